# Applies the "Updated symbol list" data refresh to the crypto price table.
# For every changed cell we:
#   1. Force the cell format to Text ("@") so Excel keeps storing the value
#      as a string (matching the original inlineStr cells) instead of
#      re-interpreting numeric-looking text (e.g. "14", "261.83", "0.45%")
#      as a number/percentage.
#   2. Write the new value as a string.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '261.83'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '0.45%'
$ws.Range('G2').NumberFormat = '@'
$ws.Range('G2').Value = '14'

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '26.69'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '-1.99%'
$ws.Range('G3').NumberFormat = '@'
$ws.Range('G3').Value = '14'

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '4.708'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '0.23%'
$ws.Range('G4').NumberFormat = '@'
$ws.Range('G4').Value = '14'

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.06075'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '-0.80%'
$ws.Range('G5').NumberFormat = '@'
$ws.Range('G5').Value = '14'

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '6.702'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '0.68%'
$ws.Range('G6').NumberFormat = '@'
$ws.Range('G6').Value = '14'

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.8509'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '-0.03%'
$ws.Range('G7').NumberFormat = '@'
$ws.Range('G7').Value = '14'

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.9093'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '-1.37%'
$ws.Range('G8').NumberFormat = '@'
$ws.Range('G8').Value = '14'

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.1403'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '-0.33%'
$ws.Range('G9').NumberFormat = '@'
$ws.Range('G9').Value = '14'

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.05092'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '10.25%'
$ws.Range('G10').NumberFormat = '@'
$ws.Range('G10').Value = '14'

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07094'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '-0.03%'
$ws.Range('G11').NumberFormat = '@'
$ws.Range('G11').Value = '14'

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.03115'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '1.18%'
$ws.Range('G12').NumberFormat = '@'
$ws.Range('G12').Value = '14'

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.09043'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '-0.12%'
$ws.Range('G13').NumberFormat = '@'
$ws.Range('G13').Value = '14'

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.001533'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '-1.20%'
$ws.Range('G14').NumberFormat = '@'
$ws.Range('G14').Value = '14'

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0006155'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '1.48%'
$ws.Range('G15').NumberFormat = '@'
$ws.Range('G15').Value = '14'

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.006032'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '-0.98%'
$ws.Range('G16').NumberFormat = '@'
$ws.Range('G16').Value = '14'

# Row 17
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '-0.03%'
$ws.Range('G17').NumberFormat = '@'
$ws.Range('G17').Value = '14'

# Row 18
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '0.19%'
$ws.Range('G18').NumberFormat = '@'
$ws.Range('G18').Value = '14'

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.147'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '-0.73%'
$ws.Range('G19').NumberFormat = '@'
$ws.Range('G19').Value = '14'

# Row 20
$ws.Range('G20').NumberFormat = '@'
$ws.Range('G20').Value = '14'

# Row 21
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '-1.82%'
$ws.Range('G21').NumberFormat = '@'
$ws.Range('G21').Value = '14'

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.103'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '0.17%'
$ws.Range('G22').NumberFormat = '@'
$ws.Range('G22').Value = '14'

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.04241'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '0.45%'
$ws.Range('G23').NumberFormat = '@'
$ws.Range('G23').Value = '14'

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.001177'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '-3.24%'
$ws.Range('G24').NumberFormat = '@'
$ws.Range('G24').Value = '14'

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.004063'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '7.14%'
$ws.Range('G25').NumberFormat = '@'
$ws.Range('G25').Value = '14'

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0001200'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '0.01%'
$ws.Range('G26').NumberFormat = '@'
$ws.Range('G26').Value = '14'

# Row 27
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '23.02%'
$ws.Range('G27').NumberFormat = '@'
$ws.Range('G27').Value = '14'

# Row 28
$ws.Range('G28').NumberFormat = '@'
$ws.Range('G28').Value = '14'

# Row 29
$ws.Range('G29').NumberFormat = '@'
$ws.Range('G29').Value = '14'

# Row 30
$ws.Range('G30').NumberFormat = '@'
$ws.Range('G30').Value = '14'

# Row 31
$ws.Range('G31').NumberFormat = '@'
$ws.Range('G31').Value = '14'

# Row 32
$ws.Range('G32').NumberFormat = '@'
$ws.Range('G32').Value = '14'

# Row 33
$ws.Range('G33').NumberFormat = '@'
$ws.Range('G33').Value = '14'

# Row 34
$ws.Range('G34').NumberFormat = '@'
$ws.Range('G34').Value = '14'

# Row 35
$ws.Range('G35').NumberFormat = '@'
$ws.Range('G35').Value = '14'

# Row 36
$ws.Range('G36').NumberFormat = '@'
$ws.Range('G36').Value = '14'

# Row 37
$ws.Range('G37').NumberFormat = '@'
$ws.Range('G37').Value = '14'

# Row 38
$ws.Range('G38').NumberFormat = '@'
$ws.Range('G38').Value = '14'

# Row 39
$ws.Range('G39').NumberFormat = '@'
$ws.Range('G39').Value = '14'

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.03954'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '2.19%'
$ws.Range('G40').NumberFormat = '@'
$ws.Range('G40').Value = '14'

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.1114'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '0.29%'
$ws.Range('G41').NumberFormat = '@'
$ws.Range('G41').Value = '14'

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.004181'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '2.55%'
$ws.Range('G42').NumberFormat = '@'
$ws.Range('G42').Value = '14'

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.002060'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '-7.07%'
$ws.Range('G43').NumberFormat = '@'
$ws.Range('G43').Value = '14'

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.01150'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '-29.61%'
$ws.Range('G44').NumberFormat = '@'
$ws.Range('G44').Value = '14'

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00005108'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '-0.98%'
$ws.Range('G45').NumberFormat = '@'
$ws.Range('G45').Value = '14'

# Row 46
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '0.01%'
$ws.Range('G46').NumberFormat = '@'
$ws.Range('G46').Value = '14'

# Row 47
$ws.Range('G47').NumberFormat = '@'
$ws.Range('G47').Value = '14'

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.2587'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '90.94%'
$ws.Range('G48').NumberFormat = '@'
$ws.Range('G48').Value = '14'

# Row 49
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '0.01%'
$ws.Range('G49').NumberFormat = '@'
$ws.Range('G49').Value = '14'

# Row 50
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '0.01%'
$ws.Range('G50').NumberFormat = '@'
$ws.Range('G50').Value = '14'

# Row 51
$ws.Range('G51').NumberFormat = '@'
$ws.Range('G51').Value = '14'
